# chore: update Sheets via scheduled runner
# Refreshes cached market-price / profit figures (columns H-N) across the
# Leve-profit worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 369002
$ws.Range("J32").Value = 369002
$ws.Range("L32").Value = 369002
$ws.Range("N32").Value = -369654
$ws.Range("H40").Value = 8529.111000000001
$ws.Range("J40").Value = 8792.666999999999
$ws.Range("L40").Value = 8792.666999999999
$ws.Range("N40").Value = -9142.666999999999
$ws.Range("H92").Value = 7143126
$ws.Range("I92").Value = 7692559
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 7692559
$ws.Range("L92").Value = 500
$ws.Range("M92").Value = -7691311
$ws.Range("N92").Value = -2996
$ws.Range("H132").Value = 1478.9111
$ws.Range("I132").Value = 1086.8536
$ws.Range("K132").Value = 3260.5608
$ws.Range("M132").Value = -730.5607999999997
$ws.Range("H138").Value = 1618.0625
$ws.Range("J138").Value = 2495.8333
$ws.Range("L138").Value = 7487.499899999999
$ws.Range("N138").Value = -17767.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 70053.07000000001
$ws.Range("I74").Value = 102029.6
$ws.Range("K74").Value = 102029.6
$ws.Range("M74").Value = -101155.6
$ws.Range("H77").Value = 70053.07000000001
$ws.Range("I77").Value = 102029.6
$ws.Range("K77").Value = 510148
$ws.Range("M77").Value = -505780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 58826.668
$ws.Range("J124").Value = 58826.668
$ws.Range("L124").Value = 58826.668
$ws.Range("N124").Value = -68646.66800000001
$ws.Range("H134").Value = 3476.8235
$ws.Range("I134").Value = 2811.1428
$ws.Range("J134").Value = 6583.3335
$ws.Range("K134").Value = 8433.428400000001
$ws.Range("L134").Value = 19750.0005
$ws.Range("M134").Value = -5898.428400000001
$ws.Range("N134").Value = -24820.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 10500
$ws.Range("I22").Value = 10000
$ws.Range("J22").Value = 11000
$ws.Range("K22").Value = 10000
$ws.Range("L22").Value = 11000
$ws.Range("M22").Value = -9650
$ws.Range("N22").Value = -11700
$ws.Range("H31").Value = 2271.7954
$ws.Range("I31").Value = 1722.4062
$ws.Range("J31").Value = 3736.8333
$ws.Range("K31").Value = 1722.4062
$ws.Range("L31").Value = 3736.8333
$ws.Range("M31").Value = -1427.4062
$ws.Range("N31").Value = -4326.8333
$ws.Range("H34").Value = 2271.7954
$ws.Range("I34").Value = 1722.4062
$ws.Range("J34").Value = 3736.8333
$ws.Range("K34").Value = 1722.4062
$ws.Range("L34").Value = 3736.8333
$ws.Range("M34").Value = -1520.4062
$ws.Range("N34").Value = -4140.8333
$ws.Range("H58").Value = 1533.5
$ws.Range("I58").Value = 1236.5385
$ws.Range("J58").Value = 1962.4445
$ws.Range("K58").Value = 1236.5385
$ws.Range("L58").Value = 1962.4445
$ws.Range("M58").Value = -1033.5385
$ws.Range("N58").Value = -2368.4445
$ws.Range("H86").Value = 5107824.5
$ws.Range("I86").Value = 5956628.5
$ws.Range("K86").Value = 5956628.5
$ws.Range("M86").Value = -5955505.5
$ws.Range("H89").Value = 5107824.5
$ws.Range("I89").Value = 5956628.5
$ws.Range("K89").Value = 29783142.5
$ws.Range("M89").Value = -29777526.5
$ws.Range("H122").Value = 2691.9375
$ws.Range("I122").Value = 2216.1428
$ws.Range("J122").Value = 6022.5
$ws.Range("K122").Value = 6648.428400000001
$ws.Range("L122").Value = 18067.5
$ws.Range("M122").Value = -4198.428400000001
$ws.Range("N122").Value = -22967.5
$ws.Range("H132").Value = 5052471.5
$ws.Range("I132").Value = 7577107.5
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 22731322.5
$ws.Range("L132").Value = 9600
$ws.Range("M132").Value = -22728792.5
$ws.Range("N132").Value = -14660
$ws.Range("H136").Value = 1533.5
$ws.Range("I136").Value = 1236.5385
$ws.Range("J136").Value = 1962.4445
$ws.Range("K136").Value = 3709.6155
$ws.Range("L136").Value = 5887.333500000001
$ws.Range("M136").Value = -1159.6155
$ws.Range("N136").Value = -10987.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 150.66667
$ws.Range("I11").Value = 100
$ws.Range("J11").Value = 252
$ws.Range("K11").Value = 300
$ws.Range("L11").Value = 756
$ws.Range("M11").Value = -160
$ws.Range("N11").Value = -1036
$ws.Range("H18").Value = 10070.363
$ws.Range("I18").Value = 11863.777
$ws.Range("K18").Value = 35591.331
$ws.Range("M18").Value = -35422.331
$ws.Range("H97").Value = 124.64286
$ws.Range("I97").Value = 120.75
$ws.Range("J97").Value = 129.83333
$ws.Range("K97").Value = 362.25
$ws.Range("L97").Value = 389.49999
$ws.Range("M97").Value = 133.75
$ws.Range("N97").Value = -1381.49999
$ws.Range("H131").Value = 1440.2307
$ws.Range("I131").Value = 1028.1111
$ws.Range("J131").Value = 2367.5
$ws.Range("K131").Value = 3084.3333
$ws.Range("L131").Value = 7102.5
$ws.Range("M131").Value = 1955.6667
$ws.Range("N131").Value = -17182.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 8449.75
$ws.Range("J18").Value = 9999.5
$ws.Range("L18").Value = 9999.5
$ws.Range("N18").Value = -10585.5
$ws.Range("H22").Value = 4478.4287
$ws.Range("J22").Value = 2987.25
$ws.Range("L22").Value = 2987.25
$ws.Range("N22").Value = -4045.25
$ws.Range("H70").Value = 95308.39999999999
$ws.Range("I70").Value = 6084.6665
$ws.Range("K70").Value = 6084.6665
$ws.Range("M70").Value = -5814.6665
$ws.Range("H73").Value = 95308.39999999999
$ws.Range("I73").Value = 6084.6665
$ws.Range("K73").Value = 6084.6665
$ws.Range("M73").Value = -5148.6665
$ws.Range("H113").Value = 2779217.8
$ws.Range("I113").Value = 1571
$ws.Range("J113").Value = 33333332
$ws.Range("K113").Value = 1571
$ws.Range("L113").Value = 33333332
$ws.Range("M113").Value = 599
$ws.Range("N113").Value = -33337672

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2855.5881
$ws.Range("I7").Value = 1996.5238
$ws.Range("J7").Value = 4243.3076
$ws.Range("K7").Value = 1996.5238
$ws.Range("L7").Value = 4243.3076
$ws.Range("M7").Value = -1884.5238
$ws.Range("N7").Value = -4467.3076
$ws.Range("H25").Value = 761.4
$ws.Range("I25").Value = 307
$ws.Range("J25").Value = 875
$ws.Range("K25").Value = 307
$ws.Range("L25").Value = 875
$ws.Range("M25").Value = -77
$ws.Range("N25").Value = -1335
$ws.Range("H46").Value = 5297.923
$ws.Range("I46").Value = 8873.154
$ws.Range("J46").Value = 1722.6923
$ws.Range("K46").Value = 8873.154
$ws.Range("L46").Value = 1722.6923
$ws.Range("M46").Value = -8685.154
$ws.Range("N46").Value = -2098.6923
$ws.Range("H93").Value = 2308.1428
$ws.Range("I93").Value = 1197.8182
$ws.Range("J93").Value = 3529.5
$ws.Range("K93").Value = 1197.8182
$ws.Range("L93").Value = 3529.5
$ws.Range("M93").Value = 50.18180000000007
$ws.Range("N93").Value = -6025.5
$ws.Range("H126").Value = 2855.5881
$ws.Range("I126").Value = 1996.5238
$ws.Range("J126").Value = 4243.3076
$ws.Range("K126").Value = 5989.5714
$ws.Range("L126").Value = 12729.9228
$ws.Range("M126").Value = -3519.5714
$ws.Range("N126").Value = -17669.9228
$ws.Range("H140").Value = 64999
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10110.286
$ws.Range("J41").Value = 10091.833
$ws.Range("L41").Value = 10091.833
$ws.Range("N41").Value = -10871.833
$ws.Range("H126").Value = 3110.4783
$ws.Range("I126").Value = 2514.75
$ws.Range("J126").Value = 4472.143
$ws.Range("K126").Value = 7544.25
$ws.Range("L126").Value = 13416.429
$ws.Range("M126").Value = -5074.25
$ws.Range("N126").Value = -18356.429
$ws.Range("H136").Value = 2617.913
$ws.Range("I136").Value = 2452.8823
$ws.Range("J136").Value = 3085.5
$ws.Range("K136").Value = 7358.646900000001
$ws.Range("L136").Value = 9256.5
$ws.Range("M136").Value = -4808.646900000001
$ws.Range("N136").Value = -14356.5
